# Hungary_FX.xlsx update:
#  - Row 313 (last existing data row): update E/F/G (low/close/volume)
#  - Append three new monthly rows: 314, 315, 316

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 313 ---
$ws.Cells.Item(313, 5).Value = 345.64
$ws.Cells.Item(313, 6).Value = 351.752
$ws.Cells.Item(313, 7).Value = 2204228

# --- New row 314 ---
# Copy formatting (date style) from A313 down to the new date cells first.
$ws.Cells.Item(313, 1).Copy()
$ws.Cells.Item(314, 1).PasteSpecial(-4122)
$ws.Cells.Item(314, 1).Value = 45170.29166666666
$ws.Cells.Item(314, 2).Value = "OANDA:USDHUF"
$ws.Cells.Item(314, 3).Value = 351.752
$ws.Cells.Item(314, 4).Value = 375.88
$ws.Cells.Item(314, 5).Value = 350.544
$ws.Cells.Item(314, 6).Value = 368.605
$ws.Cells.Item(314, 7).Value = 1642403

# --- New row 315 ---
$ws.Cells.Item(313, 1).Copy()
$ws.Cells.Item(315, 1).PasteSpecial(-4122)
$ws.Cells.Item(315, 1).Value = 45201.33333333334
$ws.Cells.Item(315, 2).Value = "OANDA:USDHUF"
$ws.Cells.Item(315, 3).Value = 368.605
$ws.Cells.Item(315, 4).Value = 375.607
$ws.Cells.Item(315, 5).Value = 356.362
$ws.Cells.Item(315, 6).Value = 361.94
$ws.Cells.Item(315, 7).Value = 2470113

# --- New row 316 ---
$ws.Cells.Item(313, 1).Copy()
$ws.Cells.Item(316, 1).PasteSpecial(-4122)
$ws.Cells.Item(316, 1).Value = 45231.33333333334
$ws.Cells.Item(316, 2).Value = "OANDA:USDHUF"
$ws.Cells.Item(316, 3).Value = 361.94
$ws.Cells.Item(316, 4).Value = 365.205
$ws.Cells.Item(316, 5).Value = 351.586
$ws.Cells.Item(316, 6).Value = 352.645
$ws.Cells.Item(316, 7).Value = 748676

$excel.CutCopyMode = $false
